$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- D2: value 2 -> 3 ---
$ws.Range("D2").Value = 3

# --- Row 3: drop the stray empty E3 cell entirely ---
$ws.Range("E3").Clear()

# --- Re-style existing rows 18-20 to match the centered/bordered-row look used by rows 7-17 ---
# (horizontal alignment = center, row height 12.75 -> matches the style already used just above them)
$ws.Rows.Item(18).RowHeight = 12.75
$ws.Range("A18:D18").HorizontalAlignment = -4108

$ws.Rows.Item(19).RowHeight = 12.75
$ws.Range("A19:D19").HorizontalAlignment = -4108

$ws.Rows.Item(20).RowHeight = 12.75
$ws.Range("A20:D20").HorizontalAlignment = -4108

# --- New rows 21-23: same centered style as rows 18-20 ---
$ws.Rows.Item(21).RowHeight = 12.75
$ws.Range("A21:D21").HorizontalAlignment = -4108
$ws.Range("A21").Value = 12.94
$ws.Range("B21").Value = 2
$ws.Range("C21").Value = "20/11/2025"
$ws.Range("D21").Value = "19:00:34"

$ws.Rows.Item(22).RowHeight = 12.75
$ws.Range("A22:D22").HorizontalAlignment = -4108
$ws.Range("A22").Value = 12.83
$ws.Range("B22").Value = 2
$ws.Range("C22").Value = "20/11/2025"
$ws.Range("D22").Value = "20:14:21"

$ws.Rows.Item(23).RowHeight = 12.75
$ws.Range("A23:D23").HorizontalAlignment = -4108
$ws.Range("A23").Value = 12.78
$ws.Range("B23").Value = 2
$ws.Range("C23").Value = "21/11/2025"
$ws.Range("D23").Value = "07:00:26"

# --- New rows 24-25: plain row (no custom height), explicit "general" alignment override ---
$ws.Range("A24:D24").Style = "Normal"
$ws.Range("A24").Value = 12.07
$ws.Range("B24").Value = 3
$ws.Range("C24").Value = "21/11/2025"
$ws.Range("D24").Value = "14:00:27"

$ws.Range("A25:D25").Style = "Normal"
$ws.Range("A25").Value = 11.92
$ws.Range("B25").Value = 3
$ws.Range("C25").Value = "22/11/2025"
$ws.Range("D25").Value = "07:00:21"

# --- New row 26: plain/general formatting, same as rows 24-25 ---
$ws.Range("A26:D26").Style = "Normal"
$ws.Range("A26").Value = 11.84
$ws.Range("B26").Value = 3
$ws.Range("C26").Value = "22/11/2025"
$ws.Range("D26").Value = "14:00:17"

# --- Selection moves to D2 ---
$ws.Range("D2").Select()
